# Developer Guide edit: refresh "today" date fields in headers/footers and
# rename two class-diagram shapes (BrowserPanel -> CalendarPanel split into
# two runs, Task/Card -> merged TaskCard run) on the UiComponentClassDiagram
# slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer/header "datetimeFigureOut" placeholders: 3/5/2018 -> 4/2/2018
#    These live on the slide master, every slide layout, and the notes
#    master. Re-stamping the whole TextRange (via a Characters() range so
#    the engine always re-materialises the run even when old == candidate
#    text) keeps the existing run formatting intact.
# ---------------------------------------------------------------------
function Update-DateText($shape) {
    if ($shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "3/5/2018") {
            $tr.Characters(1, $tr.Length).Text = "4/2/2018"
        }
    }
}

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateText $master.Shapes.Item($i)
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateText $layout.Shapes.Item($i)
    }
}

# The Notes Master's date placeholder only accepts edits through the
# HeadersFooters façade in this host.
$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = "4/2/2018"

# ---------------------------------------------------------------------
# 2) UI component class diagram shape renames (slide 1).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

# "BrowserPanel" -> split into "Calendar" + "Panel" runs.
$browserPanelShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "BrowserPanel") {
        $browserPanelShape = $shp
    }
}
if ($browserPanelShape -ne $null) {
    $tr = $browserPanelShape.TextFrame.TextRange
    # Replace the first 7 characters ("Browser") with "Calendar", leaving
    # "Panel" behind as a second, independent run.
    $tr.Characters(1, 7).Text = "Calendar"
}

# "Task" + "Card" (two runs) -> merged into a single "TaskCard" run.
$taskCardShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "TaskCard") {
        $taskCardShape = $shp
    }
}
if ($taskCardShape -ne $null) {
    $tr = $taskCardShape.TextFrame.TextRange
    # Re-stamp the whole range so the two existing runs ("Task" + "Card")
    # collapse into the single run required by the edit.
    $tr.Characters(1, $tr.Length).Text = "TaskCard"
}

Write-Host "Edit applied"
